$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters corresponding to each value in a data row (G and H are left unchanged)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New values for rows 2..25 (data rows 0..23), one array per row in the same column order as $cols
$data = @(
    @(1.02,1.0348586989063,1.03518228688624,1.052685556582033,1.058434210375272,1.031415978026885,1.039975641509913,1.037979625199116,1.055433604647065,1.061166455185883,1.041452526207731),
    @(1.02,1.036398705491944,1.036280728354912,1.054252431329086,1.060149970306888,1.031697048063639,1.041155972686871,1.038886639757722,1.056811441345373,1.062693962937536,1.042634533590346),
    @(1.02,1.037392972569845,1.036989459457073,1.055265232647773,1.061259101552275,1.031876722318253,1.04191716716419,1.039470924486024,1.0577014343427,1.063680823197225,1.043396809051133),
    @(1.02,1.037810440644078,1.037286928819531,1.055690767137729,1.061725132063025,1.031951733448762,1.042236568038801,1.039715937336161,1.058075222600292,1.064095340101512,1.043716663511621),
    @(1.02,1.037880505052463,1.037336847183029,1.055762202025601,1.061803366376928,1.031964297482231,1.042290161567373,1.039757039845187,1.058137962148778,1.064164918524633,1.043770333149151),
    @(1.02,1.037398552837078,1.036993436140966,1.055270919620216,1.061265329637529,1.031877726676875,1.041921437383724,1.039474200791013,1.05770643034256,1.063686363392021,1.043401085334868),
    @(1.02,1.035379616903307,1.035553934945533,1.053215316706337,1.05901428936664,1.031511423047827,1.040375074155712,1.038286699614192,1.055899579069422,1.061683007991642,1.041852526093702),
    @(1.02,1.031804553148846,1.033001512763185,1.049584457144414,1.055038931719709,1.030849033274505,1.037630259800056,1.036173903117914,1.052703348290705,1.05814062776088,1.039103813786826),
    @(1.02,1.029408837293565,1.031288884022143,1.047157497773954,1.052382167537013,1.030395938249248,1.035786525554226,1.034751399737719,1.050563699913741,1.055770237834187,1.037257461226887),
    @(1.02,1.028368401535679,1.030544602081528,1.046104951486897,1.051230065932639,1.030196985911631,1.034984777701611,1.034132048228731,1.049634994144926,1.054741612399846,1.036454574800361),
    @(1.02,1.027981463450201,1.030267729413838,1.045713728629044,1.050801855310405,1.030122669038401,1.03468645319292,1.033901476159345,1.049289686903293,1.054359188456337,1.03615582663665),
    @(1.02,1.028064484621766,1.030327138356414,1.045797659223051,1.050893720309309,1.030138629177091,1.034750468474347,1.033950958171809,1.049363772240048,1.054441235562819,1.036219932827119),
    @(1.02,1.02833642685011,1.030521724171918,1.046072618259682,1.051194675465555,1.030190851379779,1.034960128758114,1.034112999662813,1.049606458030501,1.054710008255877,1.036429890852537),
    @(1.02,1.02850391627418,1.030641559942816,1.046241994762301,1.051380067873935,1.03022297185277,1.035089238247028,1.034212770015696,1.049755938782699,1.05487556170193,1.036559183691734),
    @(1.02,1.02947782178547,1.031338222013185,1.047227315927514,1.052458591726586,1.030409083708095,1.035839662541792,1.034792431882884,1.050625287311882,1.055838456335216,1.037310673675071),
    @(1.02,1.03008789532949,1.031774491149485,1.047844930272213,1.053134655855689,1.030525086230553,1.036309467582625,1.035155123890688,1.051170003924254,1.056441849335459,1.037781145892949),
    @(1.02,1.030443445027365,1.032028699305496,1.048205015180762,1.053528828956184,1.0305924824923,1.036583169598945,1.035366348593622,1.051487514045056,1.056793583790278,1.038055236597527),
    @(1.02,1.030564628459217,1.03211533379118,1.048327768187347,1.053663204415698,1.030615417827111,1.036676439666748,1.03543831542919,1.05159574083081,1.056913480116963,1.038148639119525),
    @(1.02,1.030022470902364,1.031727710571051,1.047778682635755,1.053062137592967,1.030512667794877,1.036259095869064,1.035116244385433,1.051111583183397,1.056377133244064,1.037730702645778),
    @(1.02,1.028256359820786,1.030464434964256,1.045991656954053,1.05110605919483,1.030175484778003,1.034898403448302,1.034065296818699,1.049535002710176,1.054630871058925,1.036368077885704),
    @(1.02,1.027143187738467,1.029667769083383,1.044866571244518,1.049874633614156,1.02996106982565,1.034039870994295,1.033401527764081,1.048541745180967,1.053530915298023,1.035508326217138),
    @(1.02,1.027733565809632,1.030090326125929,1.045463147585494,1.050527587858421,1.03007496499369,1.03449528375067,1.033753690633622,1.0490684827676,1.054114217177234,1.035964385711863),
    @(1.02,1.03005203429914,1.031748849505127,1.047808617582731,1.053094905980175,1.030518279973544,1.036281857689684,1.035133813372099,1.051137981662387,1.056406376322874,1.037753496790794),
    @(1.02,1.032730922495528,1.033663289298952,1.050524203637959,1.056067759090974,1.031022294673981,1.038342268553124,1.036722548832691,1.053531166732691,1.059057925841902,1.039816833673993)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colLetter = $cols[$j]
        $ws.Range("$colLetter$rowNum").Value = $rowVals[$j]
    }
}
